$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Sector'
$ws.Cells.Item(1, 2).Value = 'Average Correlation'
$ws.Cells.Item(2, 1).Value = 'Containers & Packaging(12)'
$ws.Cells.Item(2, 2).Value = 0.7001438813779158
$ws.Cells.Item(3, 1).Value = 'Road & Rail(22)'
$ws.Cells.Item(3, 2).Value = 0.6293722548865807
$ws.Cells.Item(4, 1).Value = 'Energy Equipment & Services(32)'
$ws.Cells.Item(4, 2).Value = 0.6276908157365403
$ws.Cells.Item(5, 1).Value = 'Multi-Utilities(18)'
$ws.Cells.Item(5, 2).Value = 0.5813759527450517
$ws.Cells.Item(6, 1).Value = 'Building Products(23)'
$ws.Cells.Item(6, 2).Value = 0.5665676050609382
$ws.Cells.Item(7, 1).Value = 'Airlines(14)'
$ws.Cells.Item(7, 2).Value = 0.5044110339867112
$ws.Cells.Item(8, 1).Value = 'Life Sciences Tools & Services(19)'
$ws.Cells.Item(8, 2).Value = 0.4932196026268766
$ws.Cells.Item(9, 1).Value = 'Machinery(85)'
$ws.Cells.Item(9, 2).Value = 0.4888186444504815
$ws.Cells.Item(10, 1).Value = 'Internet & Direct Marketing Retail(15)'
$ws.Cells.Item(10, 2).Value = 0.4756808738508365
$ws.Cells.Item(11, 1).Value = 'Trading Companies & Distributors(25)'
$ws.Cells.Item(11, 2).Value = 0.4577437672919444
$ws.Cells.Item(12, 1).Value = 'Construction & Engineering(20)'
$ws.Cells.Item(12, 2).Value = 0.4499771060762254
$ws.Cells.Item(13, 1).Value = 'Gas Utilities(12)'
$ws.Cells.Item(13, 2).Value = 0.4484509301746201
$ws.Cells.Item(14, 1).Value = 'Semiconductors & Semiconductor Equipment(68)'
$ws.Cells.Item(14, 2).Value = 0.4397621169210953
$ws.Cells.Item(15, 1).Value = 'Auto Components(21)'
$ws.Cells.Item(15, 2).Value = 0.4372359952101255
$ws.Cells.Item(16, 1).Value = 'Consumer Finance(15)'
$ws.Cells.Item(16, 2).Value = 0.4372342122136115
$ws.Cells.Item(17, 1).Value = 'Banks(246)'
$ws.Cells.Item(17, 2).Value = 0.4338248876960353
$ws.Cells.Item(18, 1).Value = 'Leisure Products(11)'
$ws.Cells.Item(18, 2).Value = 0.4237333247381553
$ws.Cells.Item(19, 1).Value = 'Electrical Equipment(28)'
$ws.Cells.Item(19, 2).Value = 0.4227163602566927
$ws.Cells.Item(20, 1).Value = 'Household Products(10)'
$ws.Cells.Item(20, 2).Value = 0.4055378415224902
$ws.Cells.Item(21, 1).Value = 'Real Estate Management & Development(22)'
$ws.Cells.Item(21, 2).Value = 0.39762561155322
$ws.Cells.Item(22, 1).Value = 'Chemicals(51)'
$ws.Cells.Item(22, 2).Value = 0.3698598457037373
$ws.Cells.Item(23, 1).Value = 'Air Freight & Logistics(11)'
$ws.Cells.Item(23, 2).Value = 0.3677915872076659
$ws.Cells.Item(24, 1).Value = 'Oil, Gas & Consumable Fuels(122)'
$ws.Cells.Item(24, 2).Value = 0.3534335435368156
$ws.Cells.Item(25, 1).Value = 'Health Care Providers & Services(46)'
$ws.Cells.Item(25, 2).Value = 0.3426805157079772
$ws.Cells.Item(26, 1).Value = 'IT Services(52)'
$ws.Cells.Item(26, 2).Value = 0.3391849322854577
$ws.Cells.Item(27, 1).Value = 'Commercial Services & Supplies(52)'
$ws.Cells.Item(27, 2).Value = 0.3308361528058588
$ws.Cells.Item(28, 1).Value = 'Software(66)'
$ws.Cells.Item(28, 2).Value = 0.3306190728028508
$ws.Cells.Item(29, 1).Value = 'Hotels, Restaurants & Leisure(50)'
$ws.Cells.Item(29, 2).Value = 0.326760669788997
$ws.Cells.Item(30, 1).Value = 'Electric Utilities(28)'
$ws.Cells.Item(30, 2).Value = 0.324950483858519
$ws.Cells.Item(31, 1).Value = 'Capital Markets(75)'
$ws.Cells.Item(31, 2).Value = 0.3211408690304468
$ws.Cells.Item(32, 1).Value = 'Household Durables(39)'
$ws.Cells.Item(32, 2).Value = 0.3181525400969379
$ws.Cells.Item(33, 1).Value = 'Beverages(21)'
$ws.Cells.Item(33, 2).Value = 0.3127055360960744
$ws.Cells.Item(34, 1).Value = 'Aerospace & Defense(37)'
$ws.Cells.Item(34, 2).Value = 0.3087510583339029
$ws.Cells.Item(35, 1).Value = 'Insurance(75)'
$ws.Cells.Item(35, 2).Value = 0.3024132463637215
$ws.Cells.Item(36, 1).Value = 'Professional Services(35)'
$ws.Cells.Item(36, 2).Value = 0.2985858803840739
$ws.Cells.Item(37, 1).Value = 'Thrifts & Mortgage Finance(47)'
$ws.Cells.Item(37, 2).Value = 0.2939761905456659
$ws.Cells.Item(38, 1).Value = 'Marine(15)'
$ws.Cells.Item(38, 2).Value = 0.2834946746149619
$ws.Cells.Item(39, 1).Value = 'Textiles, Apparel & Luxury Goods(29)'
$ws.Cells.Item(39, 2).Value = 0.2794280827201666
$ws.Cells.Item(40, 1).Value = 'Health Care Equipment & Supplies(83)'
$ws.Cells.Item(40, 2).Value = 0.2730344860293868
$ws.Cells.Item(41, 1).Value = 'Food Products(44)'
$ws.Cells.Item(41, 2).Value = 0.2677485921181879
$ws.Cells.Item(42, 1).Value = 'Biotechnology(126)'
$ws.Cells.Item(42, 2).Value = 0.24546180325089
$ws.Cells.Item(43, 1).Value = 'Metals & Mining(89)'
$ws.Cells.Item(43, 2).Value = 0.243965935247877
$ws.Cells.Item(44, 1).Value = 'Media(42)'
$ws.Cells.Item(44, 2).Value = 0.2372918794920321
$ws.Cells.Item(45, 1).Value = 'Communications Equipment(45)'
$ws.Cells.Item(45, 2).Value = 0.2308892914139308
$ws.Cells.Item(46, 1).Value = 'Specialty Retail(58)'
$ws.Cells.Item(46, 2).Value = 0.2240830577663852
$ws.Cells.Item(47, 1).Value = 'Entertainment(22)'
$ws.Cells.Item(47, 2).Value = 0.2067035958370407
$ws.Cells.Item(48, 1).Value = 'Pharmaceuticals(48)'
$ws.Cells.Item(48, 2).Value = 0.1815443032280717

# Remove now-unused trailing rows (49:53) that existed before the data refresh
$ws.Rows("49:53").Delete()
